$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.306533813476562
$ws.Range("B1").Value = 2.525996923446655
$ws.Range("C1").Value = 2.561086416244507
$ws.Range("D1").Value = 3.251309633255005
$ws.Range("E1").Value = 2.301841020584106
